{"js": "// \"Try to fix Error\" \u2014 Custom ErrorCodeDetails.docx\n//\n// 1. \"Sql Grammar Exception\"            -> \"Sql\" flagged spellStart/spellEnd,\n//                                           then \" Grammar Exception\"\n// 2. \"Network or driver issue or db is temporarily unavailable\"\n//                                       -> \"Network or driver issue or \" +\n//                                          \"db\" flagged spellStart/spellEnd +\n//                                          \" is temporarily unavailable\"\n// 3. Append a new table row: 4005 / \"Operation Failed. You are not permitted\n//    to update vehicle settings\"\n\nconst WORD_NS = 'xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\"';\n\n// Office.js can only inject raw WordprocessingML via Range.insertOoxml, and\n// that API requires the FlatOpc <pkg:package> envelope (InsertHtml is\n// membrane-blocked, and there's no direct \"add a proofErr\" call).\nfunction wrapFlatOpc(bodyInnerXml) {\n  return (\n    '<?xml version=\"1.0\" encoding=\"UTF-8\" standalone=\"yes\"?>' +\n    '<pkg:package xmlns:pkg=\"http://schemas.microsoft.com/office/2006/xmlPackage\">' +\n    '<pkg:part pkg:name=\"/word/document.xml\" ' +\n    'pkg:contentType=\"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml\">' +\n    \"<pkg:xmlData>\" +\n    `<w:document ${WORD_NS}><w:body>${bodyInnerXml}</w:body></w:document>` +\n    \"</pkg:xmlData>\" +\n    \"</pkg:part>\" +\n    \"</pkg:package>\"\n  );\n}\n\nasync function replaceWithOoxmlParagraph(body, searchText, paragraphInnerXml) {\n  const results = body.search(searchText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n  if (results.items.length === 0) {\n    throw new Error(`Text not found: ${searchText}`);\n  }\n  const range = results.items[0];\n  const flatOpc = wrapFlatOpc(`<w:p>${paragraphInnerXml}</w:p>`);\n  range.insertOoxml(flatOpc, Word.InsertLocation.replace);\n  await context.sync();\n}\n\nconst body = context.document.body;\n\n// --- 1. Split \"Sql Grammar Exception\", flag \"Sql\" as a spell-check error ---\nawait replaceWithOoxmlParagraph(\n  body,\n  \"Sql Grammar Exception\",\n  '<w:proofErr w:type=\"spellStart\"/>' +\n    \"<w:r><w:t>Sql</w:t></w:r>\" +\n    '<w:proofErr w:type=\"spellEnd\"/>' +\n    '<w:r><w:t xml:space=\"preserve\"> Grammar Exception</w:t></w:r>'\n);\n\n// --- 2. Split the \"db\" message, flag \"db\" as a spell-check error ---\nawait replaceWithOoxmlParagraph(\n  body,\n  \"Network or driver issue or db is temporarily unavailable\",\n  '<w:r><w:t xml:space=\"preserve\">Network or driver issue or </w:t></w:r>' +\n    '<w:proofErr w:type=\"spellStart\"/>' +\n    \"<w:r><w:t>db</w:t></w:r>\" +\n    '<w:proofErr w:type=\"spellEnd\"/>' +\n    '<w:r><w:t xml:space=\"preserve\"> is temporarily unavailable</w:t></w:r>'\n);\n\n// --- 3. Append a new row with the new error code / message ---\nconst tables = body.tables;\ntables.load(\"items\");\nawait context.sync();\nconst table = tables.items[0];\ntable.addRows(\"End\", 1, [\n  [\"4005\", \"Operation Failed. You are not permitted to update vehicle settings\"],\n]);\nawait context.sync();\n", "ps1": "# \"Try to fix Error\" \u2014 Custom ErrorCodeDetails.docx\n#\n# 1. \"Sql Grammar Exception\"            -> \"Sql\" flagged spellStart/spellEnd,\n#                                           then \" Grammar Exception\"\n# 2. \"Network or driver issue or db is temporarily unavailable\"\n#                                       -> \"Network or driver issue or \" +\n#                                          \"db\" flagged spellStart/spellEnd +\n#                                          \" is temporarily unavailable\"\n# 3. Append a new table row: 4005 / \"Operation Failed. You are not permitted\n#    to update vehicle settings\"\n\n$d = $word.ActiveDocument\n$wordNs = 'http://schemas.openxmlformats.org/wordprocessingml/2006/main'\n\n# --- 1. Split \"Sql Grammar Exception\", flag \"Sql\" as a spell-check error ---\n$rng = $d.Content\n$find = $rng.Find\n$find.ClearFormatting()\n$find.Text = \"Sql Grammar Exception\"\n$found = $find.Execute()\nif ($found) {\n    $xml = '<w:p xmlns:w=\"' + $wordNs + '\">' +\n        '<w:proofErr w:type=\"spellStart\"/>' +\n        '<w:r><w:t>Sql</w:t></w:r>' +\n        '<w:proofErr w:type=\"spellEnd\"/>' +\n        '<w:r><w:t xml:space=\"preserve\"> Grammar Exception</w:t></w:r>' +\n        '</w:p>'\n    $rng.InsertXML($xml)\n}\n\n# --- 2. Split the \"db\" message, flag \"db\" as a spell-check error ---\n$rng2 = $d.Content\n$find2 = $rng2.Find\n$find2.ClearFormatting()\n$find2.Text = \"Network or driver issue or db is temporarily unavailable\"\n$found2 = $find2.Execute()\nif ($found2) {\n    $xml2 = '<w:p xmlns:w=\"' + $wordNs + '\">' +\n        '<w:r><w:t xml:space=\"preserve\">Network or driver issue or </w:t></w:r>' +\n        '<w:proofErr w:type=\"spellStart\"/>' +\n        '<w:r><w:t>db</w:t></w:r>' +\n        '<w:proofErr w:type=\"spellEnd\"/>' +\n        '<w:r><w:t xml:space=\"preserve\"> is temporarily unavailable</w:t></w:r>' +\n        '</w:p>'\n    $rng2.InsertXML($xml2)\n}\n\n# --- 3. Append a new row with the new error code / message ---\n$table = $d.Tables.Item(1)\n$newRow = $table.Rows.Add()\n$rowIndex = $table.Rows.Count\n$table.Cell($rowIndex, 1).Range.Text = \"4005\"\n$table.Cell($rowIndex, 2).Range.Text = \"Operation Failed. You are not permitted to update vehicle settings\"\n"}
